# Fruta / hortaliza, semanal
# Inserts two new weekly price records for "Feria Lagunitas de Puerto Montt - Pomelo"
# right before the existing row 307, pushing the remainder of the table down by
# two rows (old row 307 -> new row 309, ..., old row 321 -> new row 323).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 307 (Excel shifts 307..321 down to 309..323).
$ws.Rows("307:308").Insert()

# New row 307: Start Ruby / Primera
$ws.Cells.Item(307, 1).Value = 4
$ws.Cells.Item(307, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(307, 3).Value = "Los Lagos"
$ws.Cells.Item(307, 4).Value = 44753
$ws.Cells.Item(307, 5).Value = 10
$ws.Cells.Item(307, 6).Value = "Fruta"
$ws.Cells.Item(307, 7).Value = 100102
$ws.Cells.Item(307, 8).Value = "Cítricos"
$ws.Cells.Item(307, 9).Value = 100102006
$ws.Cells.Item(307, 10).Value = "Pomelo"
$ws.Cells.Item(307, 11).Value = "Start Ruby"
$ws.Cells.Item(307, 12).Value = "Primera"
$ws.Cells.Item(307, 13).Value = 50
$ws.Cells.Item(307, 14).Value = 15000
$ws.Cells.Item(307, 15).Value = 15000
$ws.Cells.Item(307, 16).Value = 15000
$ws.Cells.Item(307, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(307, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(307, 19).Value = 1071
$ws.Cells.Item(307, 20).Value = 14

# New row 308: Start Ruby / Segunda
$ws.Cells.Item(308, 1).Value = 4
$ws.Cells.Item(308, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(308, 3).Value = "Los Lagos"
$ws.Cells.Item(308, 4).Value = 44753
$ws.Cells.Item(308, 5).Value = 10
$ws.Cells.Item(308, 6).Value = "Fruta"
$ws.Cells.Item(308, 7).Value = 100102
$ws.Cells.Item(308, 8).Value = "Cítricos"
$ws.Cells.Item(308, 9).Value = 100102006
$ws.Cells.Item(308, 10).Value = "Pomelo"
$ws.Cells.Item(308, 11).Value = "Start Ruby"
$ws.Cells.Item(308, 12).Value = "Segunda"
$ws.Cells.Item(308, 13).Value = 50
$ws.Cells.Item(308, 14).Value = 13000
$ws.Cells.Item(308, 15).Value = 13000
$ws.Cells.Item(308, 16).Value = 13000
$ws.Cells.Item(308, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(308, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(308, 19).Value = 929
$ws.Cells.Item(308, 20).Value = 14
